# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" right before the "总计" sheet, cloned from
#    "2021-Q4" so it inherits the same layout/styles, then overwrite its
#    contents with the new quarter's fund-holding table.
# 2. Extend "总计" with a new top data row for 2022-Q1 and shift the
#    existing history rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q1" sheet
# ---------------------------------------------------------------------

$sourceSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# Clone the most recent quarterly sheet so the new tab starts out with the
# identical column layout / fonts / borders, then drop it right before 总计.
$sourceSheet.Copy($totalSheet)
$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"

# The template sheet only has 10 data rows (rows 2-11); we need 14 data
# rows (rows 2-15). Extend the formatting of the last template row down.
$ws.Range("A10:H10").Copy()
$ws.Range("A11:H15").PasteSpecial(-4122)

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$headerCols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $headers[$i]
}

$fundData = @(
    @("501025", "鹏华港股通中证香港银行投资指数（LOF）A", "9.81", "94.47", "4.42", "0.4336", 8),
    @("010365", "鹏华港股通中证香港银行投资指数（LOF）C", "6.07", "94.47", "4.42", "0.2683", 8),
    @("010204", "中银港股通优势成长股票", "3.19", "83.00", "2.86", "0.0912", 10),
    @("006809", "泰康港股通中证香港银行投资指数A", "1.99", "94.73", "4.42", "0.0880", 8),
    @("501310", "华宝标普沪港深中国增强价值指数（LOF）A", "1.44", "94.80", "3.83", "0.0552", 5),
    @("501305", "汇添富中证港股通高股息投资指数（LOF）A", "1.59", "93.08", "3.06", "0.0487", 10),
    @("006810", "泰康港股通中证香港银行投资指数C", "0.90", "94.73", "4.42", "0.0398", 8),
    @("007751", "景顺长城中证沪港深红利成长低波动指数A", "0.83", "91.29", "3.41", "0.0283", 2),
    @("006658", "财通中证香港红利等权投资指数A", "0.20", "90.59", "3.48", "0.0070", 6),
    @("501306", "汇添富中证港股通高股息投资指数（LOF）C", "0.21", "93.08", "3.06", "0.0064", 10),
    @("007397", "华宝标普沪港深中国增强价值指数（LOF）C", "0.09", "94.80", "3.83", "0.0034", 5),
    @("005269", "华泰柏瑞港股通量化灵活配置混合", "0.33", "37.77", "0.86", "0.0028", 9),
    @("007760", "景顺长城中证沪港深红利成长低波动指数C", "0.06", "91.29", "3.41", "0.0020", 2),
    @("006659", "财通中证香港红利等权投资指数C", "0.05", "90.59", "3.48", "0.0017", 6)
)

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $row = 2 + $i
    $rec = $fundData[$i]

    # Column A: numeric zero-based rank, keeps the s="2" style already on it.
    $ws.Range("A" + $row).Value = $i

    # Columns B, C, D, E, F, G hold text even when the text looks numeric
    # (fund code / name / size / position figures) - force text storage
    # with a leading apostrophe, then strip the Text number-format so the
    # cell keeps the default (no explicit) style, matching the template.
    $textCols = @("B", "C", "D", "E", "F", "G")
    for ($c = 0; $c -lt $textCols.Length; $c++) {
        $cell = $ws.Range($textCols[$c] + $row)
        $cell.Value = "'" + $rec[$c]
        $cell.Style = "Normal"
    }

    # Column H is a real number.
    $ws.Range("H" + $row).Value = $rec[6]
}

# ---------------------------------------------------------------------
# Part 2: "总计" sheet - add the 2022-Q1 summary row on top
# ---------------------------------------------------------------------

$tw = $wb.Worksheets.Item("总计")

# Grow the table by one row, copying the format of the last existing row
# (keeps column A's bold/centered/bordered style).
$tw.Range("A6:D6").Copy()
$tw.Range("A7:D7").PasteSpecial(-4122)

# Shift the five existing history rows down by one (bottom-up so the
# source cells aren't clobbered before being read), carrying values +
# formats together.
for ($r = 6; $r -ge 2; $r--) {
    $destRow = $r + 1
    $tw.Range("A" + $r + ":D" + $r).Copy()
    $tw.Range("A" + $destRow + ":D" + $destRow).PasteSpecial(-4122)
    $tw.Range("A" + $r + ":D" + $r).Copy()
    $tw.Range("A" + $destRow + ":D" + $destRow).PasteSpecial(-4163)
}

# Column A is a zero-based row counter - bump every carried-down row by one.
for ($r = 7; $r -ge 3; $r--) {
    $tw.Range("A" + $r).Value = $tw.Range("A" + ($r - 1)).Value
}

# Fill in the brand-new top row.
$tw.Range("A2").Value = 0
$tw.Range("B2").Value = "2022-Q1"
$tw.Range("C2").Value = 14
$tw.Range("D2").Value = 1.08

# ---------------------------------------------------------------------
# Restore the original active sheet/selection (2020-Q4 was active before).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Select()
